# Rotate the three worker data rows (16, 17, 18) on the active sheet.
# Old row 18 -> new row 16
# Old row 16 -> new row 17
# Old row 17 -> new row 18
# Only the data values move; cell styles/number formats stay fixed per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current ("before") values for columns B..G on rows 16, 17 and 18.
$row16 = @($ws.Range("B16").Value(), $ws.Range("C16").Value(), $ws.Range("D16").Value(), $ws.Range("E16").Value(), $ws.Range("F16").Value(), $ws.Range("G16").Value())
$row17 = @($ws.Range("B17").Value(), $ws.Range("C17").Value(), $ws.Range("D17").Value(), $ws.Range("E17").Value(), $ws.Range("F17").Value(), $ws.Range("G17").Value())
$row18 = @($ws.Range("B18").Value(), $ws.Range("C18").Value(), $ws.Range("D18").Value(), $ws.Range("E18").Value(), $ws.Range("F18").Value(), $ws.Range("G18").Value())

$cols = @("B", "C", "D", "E", "F", "G")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $col = $cols[$i]
    $ws.Range($col + "16").Value = $row18[$i]
    $ws.Range($col + "17").Value = $row16[$i]
    $ws.Range($col + "18").Value = $row17[$i]
}
